# cleaning recordings 3 September 24
#
# Applies the "amendments" worksheet review pass:
#   - mark the id column (B) of several rows with a yellow highlight
#     (already-reviewed / renumbered rows) and one row with a red
#     highlight (needs a closer look)
#   - fill in the new "smp.notes" review column (E) for rows 38-48
#   - move the active selection down to where work left off (B49),
#     scrolled so row 46 is at the top of the view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ammendments")

# --- Highlight the id (column B) cells that were reviewed/renumbered ---
# Yellow fill (matches the sheet's existing "reviewed" style)
$yellow = 65535   # RGB(255,255,0)
# Red fill (new "needs re-check" style)
$red = 255        # RGB(255,0,0)

$ws.Range("B34").Interior.Color = $yellow
$ws.Range("B35").Interior.Color = $yellow
$ws.Range("B36").Interior.Color = $yellow
$ws.Range("B37").Interior.Color = $yellow
$ws.Range("B44").Interior.Color = $yellow
$ws.Range("B47").Interior.Color = $yellow
$ws.Range("B46").Interior.Color = $red

# --- Fill in the smp.notes column (E) for rows 38-48 ---
$ws.Range("E38").Value = "done"
$ws.Range("E39").Value = "done"
$ws.Range("E40").Value = "done"
$ws.Range("E41").Value = "done; aliasing"
$ws.Range("E42").Value = "done"
$ws.Range("E43").Value = "done"
$ws.Range("E44").Value = "timer"
$ws.Range("E45").Value = "done; are we including 5-note songs? Also I don't see 0.18"
$ws.Range("E46").Value = "re-check this; did not clean because the sound at 3.5 kHz seems to be part of the signal"
$ws.Range("E47").Value = "done"
$ws.Range("E48").Value = "done"

# --- Move the view / selection to reflect where work stopped ---
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 46
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B49").Select()
